# expansao das analises automaticas
# Adds three new computed-metric columns (L:N) to the existing summary
# table: apoio_medio, contribuicoes, media_contribuicoes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column titles ---
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Reuse the same header formatting (bold font + border + centered/top
# alignment) already applied to the other header cells (e.g. K1).
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (2-7): values for the three new columns ---
$data = @(
    @(88.07870613099213, 225303, 307.3710777626194),
    @(110.2975973828001, 38250,  394.3298969072165),
    @(88.66083985762999, 187667, 146.2720187061574),
    @(107.8551914385913, 15979,  159.79),
    @(17.98549503340952, 1940,   14.47761194029851),
    @(30.9518559327251,  268,    14.88888888888889)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $data[$i][0]  # L: apoio_medio
    $ws.Cells.Item($row, 13).Value = $data[$i][1]  # M: contribuicoes
    $ws.Cells.Item($row, 14).Value = $data[$i][2]  # N: media_contribuicoes
}
